$d = $word.ActiveDocument

# --- Change 1 -------------------------------------------------------------
# "higher level" -> "higher-level"
# (span across the neighbouring run so the stray gramStart/gramEnd
#  proofing-error bookmarks that used to wrap "higher level" get folded
#  away along with it, the same way Word collapses them once the phrase
#  underneath is re-typed/edited)
$found1 = $d.Content.Find.Execute(
    "more akin to the higher level protocols",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "more akin to the higher-level protocols", 2)
if (-not $found1) { Write-Output "WARN: change 1 text not found" }

# --- Change 2 ---------------------------------------------------------------
# "... letting the creator know which if it accepted or rejected the job. "
#   -> "... letting the creator know if the job is accepted or rejected. "
$found2 = $d.Content.Find.Execute(
    "letting the creator know which if it accepted or rejected the job.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "letting the creator know if the job is accepted or rejected.", 2)
if (-not $found2) { Write-Output "WARN: change 2 text not found" }

# --- Change 3 ---------------------------------------------------------------
# "job seekers and job creators task" -> "job seekers and job creators' task"
# (curly right single quotation mark, U+2019, inserted right after "creators")
$apos = [string][char]0x2019
$found3 = $d.Content.Find.Execute(
    "job seekers and job creators task",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "job seekers and job creators" + $apos + " task", 2)
if (-not $found3) { Write-Output "WARN: change 3 text not found" }
